# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'90.319.31"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "'3.088.95"
$ws.Range("E3").Value = "  -2.50%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'233.48"
$ws.Range("E5").Value = "  +7.55%  "
$ws.Range("D6").Value = "'625.27"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "'1.11"
$ws.Range("E7").Value = "  -6.56%  "
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "'0.729"
$ws.Range("E10").Value = "  -5.04%  "
$ws.Range("D11").Value = "'2.498.39"
$ws.Range("E11").Value = "  -21.05%  "
$ws.Range("D12").Value = "'0.197"
$ws.Range("E12").Value = "  -3.05%  "
$ws.Range("D13").Value = "'36.48"
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("D16").Value = "'90.126.50"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("D18").Value = "'3.079.00"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("D19").Value = "'3.76"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "'0.0000213"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").Value = "'438.39"
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("E23").Value = "  +6.90%  "
$ws.Range("D24").Value = "'8.89"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "'5.89"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").Value = "'89.09"
$ws.Range("E27").Value = "  -3.32%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'3.251.53"
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'9.48"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("D33").Value = "'0.982"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("D34").Value = "'0.206"
$ws.Range("E34").Value = "  +8.24%  "
$ws.Range("D35").Value = "'26.21"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("E36").Value = "  +6.09%  "
$ws.Range("D37").Value = "'3.84"
$ws.Range("E37").Value = "  +4.43%  "
$ws.Range("D38").Value = "'504.15"
$ws.Range("E38").Value = "  -4.07%  "
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("D40").Value = "'7.01"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("E43").Value = "  -3.13%  "
$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D44").Value = "'3.52"
$ws.Range("E44").Value = "  +52.90%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "'22.19"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("D48").Value = "'150.77"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").Value = "'0.691"
$ws.Range("E49").Value = "  +5.06%  "
$ws.Range("D50").Value = "'44.98"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").Value = "  -1.18%  "
